$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Username column (B) with the new sample usernames.
$ws.Range("B2").Value = "johnsmith99"
$ws.Range("B3").Value = "amyjackson"
$ws.Range("B4").Value = "william"
$ws.Range("B5").Value = "sophiacathie"
$ws.Range("B6").Value = "jacobmason"
$ws.Range("B7").Value = "emmawilliam"
$ws.Range("B8").Value = "Oliviajaden"
$ws.Range("B9").Value = "Emilysmith"
$ws.Range("B10").Value = "Natasha99"
$ws.Range("B11").Value = "danieldanny"

# Re-create the mailto hyperlinks on column D so they come back in the same
# order Excel wrote them in after the edit above.
$order = @(
    @(13, "mailto:Dali@gmail.com"),
    @(12, "mailto:Daley@gmail.com"),
    @(14, "mailto:Gary@gmail.com"),
    @(15, "mailto:emily@gmail.com"),
    @(11, "mailto:Curie@gmail.com"),
    @(10, "mailto:Joseph@gmail.com"),
    @(9,  "mailto:thomas@gmail.com"),
    @(8,  "mailto:richard@gmail.com"),
    @(7,  "mailto:david@gmail.com"),
    @(6,  "mailto:william@gmail.com"),
    @(5,  "mailto:michael@gmail.com"),
    @(4,  "mailto:Robert@gmail.com"),
    @(3,  "mailto:James@gmail.com"),
    @(2,  "mailto:Bob@gmail.com")
)

$ws.Hyperlinks.Delete()
foreach ($pair in $order) {
    $row = $pair[0]
    $addr = $pair[1]
    $cell = $ws.Cells.Item($row, 4)
    $ws.Hyperlinks.Add($cell, $addr) | Out-Null
    # Adding the hyperlink resets the cell's style; restore the original
    # centered alignment that the workbook's style 6 already used.
    $cell.HorizontalAlignment = -4108  # xlCenter
}

# Leave the active cell on B2, matching the saved selection state.
$ws.Range("B2").Select() | Out-Null
